$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row's updated Price (D) and Volume(1h) (E) text. Both columns hold
# plain text in the source data (even values that look numeric, e.g.
# "243.97"), so the Price column is explicitly formatted as Text before
# the value is written - otherwise Excel would silently reinterpret a
# numeric-looking string as a number.
$updates = @(
    @{ Row = 2;  D = "29.174.40";     E = "  +0.28%  " }
    @{ Row = 3;  D = "1.840.78";      E = "  +0.24%  " }
    @{ Row = 4;  D = "0.9998";        E = "  -0.24%  " }
    @{ Row = 5;  D = "243.97";        E = "  -0.34%  " }
    @{ Row = 6;  D = "0.6264";        E = "  -1.46%  " }
    @{ Row = 7;  E = "  +0.05%  " }
    @{ Row = 8;  D = "0.07539";       E = "  -0.86%  " }
    @{ Row = 9;  D = "0.2947";        E = "  -0.30%  " }
    @{ Row = 10; D = "23.36";         E = "  +2.13%  " }
    @{ Row = 11; D = "0.07709";       E = "  -0.67%  " }
    @{ Row = 12; D = "1.852.87";      E = "  +1.56%  " }
    @{ Row = 13; D = "5.028";         E = "  +0.34%  " }
    @{ Row = 14; D = "0.6784";        E = "  +0.90%  " }
    @{ Row = 15; D = "83.26";         E = "  -0.16%  " }
    @{ Row = 16; D = "0.000009314";   E = "  -5.29%  " }
    @{ Row = 17; D = "5.985";         E = "  -2.26%  " }
    @{ Row = 18; D = "29.186.43";     E = "  +0.31%  " }
    @{ Row = 19; D = "2.092.56";      E = "  +0.02%  " }
    @{ Row = 20; D = "231.99";        E = "  +2.05%  " }
    @{ Row = 21; D = "12.72";         E = "  +1.04%  " }
    @{ Row = 22; E = "  +0.03%  " }
    @{ Row = 23; D = "7.185";         E = "  -0.63%  " }
    @{ Row = 24; D = "1.001";         E = "  -0.24%  " }
    @{ Row = 25; D = "160.70";        E = "  +0.02%  " }
    @{ Row = 26; E = "  -0.04%  " }
    @{ Row = 27; D = "8.566";         E = "  +0.20%  " }
    @{ Row = 28; D = "17.96";         E = "  -0.43%  " }
    @{ Row = 29; D = "1.496";         E = "  -0.26%  " }
    @{ Row = 30; D = "4.193";         E = "  +1.48%  " }
    @{ Row = 31; D = "4.155";         E = "  +2.18%  " }
    @{ Row = 32; D = "0.05580";       E = "  +3.36%  " }
    @{ Row = 33; D = "1.209";         E = "  -0.01%  " }
    @{ Row = 34; D = "0.7502";        E = "  -0.02%  " }
    @{ Row = 35; D = "1.853";         E = "  -0.74%  " }
    @{ Row = 36; D = "1.148";         E = "  +0.45%  " }
    @{ Row = 37; D = "2.671";         E = "  +0.15%  " }
    @{ Row = 38; D = "1.241.93";      E = "  +0.40%  " }
    @{ Row = 39; D = "2.773";         E = "  +0.26%  " }
    @{ Row = 40; D = "0.01795";       E = "  -0.11%  " }
    @{ Row = 41; D = "6.617";         E = "  -0.06%  " }
    @{ Row = 42; D = "0.9022";        E = "  -0.35%  " }
    @{ Row = 43; E = "  -0.12%  " }
    @{ Row = 44; D = "102.22";        E = "  -0.35%  " }
    @{ Row = 45; E = "  +2.98%  " }
    @{ Row = 46; D = "1.988.35";      E = "  +0.01%  " }
    @{ Row = 47; D = "0.00000000123"; E = "  -0.79%  " }
    @{ Row = 48; D = "0.5092";        E = "  -0.63%  " }
    @{ Row = 49; D = "0.4097";        E = "  +0.05%  " }
    @{ Row = 50; D = "9.116";         E = "  +0.34%  " }
    @{ Row = 51; D = "0.07414";       E = "  +15.01%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $u["D"]
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$row").Value = $u["E"]
    }
}
